$d = $word.ActiveDocument

# --- Edit 1: "4%CoTraoDoiBai" -> "4%TraoDoiBai" --------------------------
$d.Content.Find.Execute("4%CoTraoDoiBai", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4%TraoDoiBai", 2) | Out-Null

# --- Edit 2: "...Client3$TongDiem%Client3$TongDiem"
#          -> "...Client3$TongDiem%Client4$TongDiem"  (split into 3 runs)
#
# Only the very last "3" (right after the final "Client") turns into "4".
# In Word that single-character retype leaves the paragraph as three runs:
#   ["...Client"] ["4"] ["$TongDiem"]
# all carrying the same (color/size) run formatting as the original run.
# We reproduce that run split by editing just that one character with
# change-tracking on, then accepting only the revisions that edit created
# (instead of Document.AcceptAllRevisions, which would also renumber/strip
# unrelated rsid bookkeeping elsewhere in the file).
$outer = $d.Content
$outerFound = $outer.Find.Execute("Client3`$TongDiem%Client3`$TongDiem",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

if ($outerFound) {
    $suffix = "3`$TongDiem"
    $charStart = $outer.End - $suffix.Length
    $charRange = $d.Range($charStart, $charStart + 1)

    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true
    $charRange.Text = "4"
    $d.TrackRevisions = $wasTracking

    foreach ($rev in $d.Revisions) {
        $rev.Accept()
    }
}
